$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.78"
$ws.Range("E2").Value = "'0.57%"
$ws.Range("E3").Value = "'-1.33%"
$ws.Range("D4").Value = "'4.673"
$ws.Range("E4").Value = "'-10.42%"
$ws.Range("D5").Value = "'0.05890"
$ws.Range("E5").Value = "'-0.55%"
$ws.Range("D6").Value = "'6.644"
$ws.Range("E6").Value = "'-0.54%"
$ws.Range("D7").Value = "'0.8578"
$ws.Range("E7").Value = "'-0.99%"
$ws.Range("D8").Value = "'0.9433"
$ws.Range("E8").Value = "'-8.71%"
$ws.Range("D9").Value = "'0.1407"
$ws.Range("E9").Value = "'-0.87%"
$ws.Range("D10").Value = "'0.03808"
$ws.Range("E10").Value = "'4.83%"
$ws.Range("D11").Value = "'0.07086"
$ws.Range("E11").Value = "'-1.34%"
$ws.Range("D12").Value = "'0.03186"
$ws.Range("E12").Value = "'-1.71%"
$ws.Range("D13").Value = "'0.09163"
$ws.Range("E13").Value = "'-0.52%"
$ws.Range("D14").Value = "'0.001557"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("D15").Value = "'0.0006039"
$ws.Range("E15").Value = "'-94.22%"
$ws.Range("D16").Value = "'0.006203"
$ws.Range("E16").Value = "'8.49%"
$ws.Range("D17").Value = "'3.513"
$ws.Range("E17").Value = "'0.80%"
$ws.Range("E18").Value = "'-1.88%"
$ws.Range("E19").Value = "'0.89%"
$ws.Range("D20").Value = "'0.3081"
$ws.Range("E20").Value = "'-2.21%"
$ws.Range("D21").Value = "'0.1293"
$ws.Range("E21").Value = "'-1.24%"
$ws.Range("D22").Value = "'3.880"
$ws.Range("E22").Value = "'10.05%"
$ws.Range("D23").Value = "'0.04232"
$ws.Range("E23").Value = "'1.37%"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("E24").Value = "'0.21%"
$ws.Range("D25").Value = "'0.004297"
$ws.Range("E25").Value = "'-4.82%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.11%"
$ws.Range("E27").Value = "'-0.13%"
$ws.Range("D40").Value = "'0.03836"
$ws.Range("E40").Value = "'0.63%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006277"
$ws.Range("E41").Value = "'14.06%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("E42").Value = "'0.17%"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'15.67%"
$ws.Range("D44").Value = "'0.01146"
$ws.Range("E44").Value = "'7.32%"
$ws.Range("D45").Value = "'0.00005462"
$ws.Range("E45").Value = "'0.52%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("D47").Value = "'0.05999"
$ws.Range("E47").Value = "'-45.05%"
$ws.Range("D48").Value = "'0.1303"
$ws.Range("E48").Value = "'5,909.83%"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("E50").Value = "'-0.10%"
